$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Target cluster" value from "ECs" to "MuSCs" for both data rows
$ws.Range("D2").Value = "MuSCs"
$ws.Range("D3").Value = "MuSCs"

# Row 2 updated metric values
$ws.Range("G2").Value = 0.092262
$ws.Range("H2").Value = 0.276786
$ws.Range("I2").Value = 0.674176843971804
$ws.Range("J2").Value = 0.6741768439718039
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.08043649999999999
$ws.Range("N2").Value = 0.160873
$ws.Range("Q2").Value = 0.007421232362999999
$ws.Range("R2").Value = 0.04452739417799999
$ws.Range("S2").Value = 0.674176843971804
$ws.Range("T2").Value = 0.6741768439718039

# Row 3 updated metric values
$ws.Range("I3").Value = 0.3258231560281961
$ws.Range("J3").Value = 0.3258231560281961
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.08043649999999999
$ws.Range("N3").Value = 0.160873
$ws.Range("Q3").Value = 0.003586609910666666
$ws.Range("R3").Value = 0.021519659464
$ws.Range("S3").Value = 0.3258231560281961
$ws.Range("T3").Value = 0.3258231560281961
